# Peru Liga 1 - base update (20-06-2024 20:11)
# The upstream scraper re-ran and a handful of fixtures on the same match-day
# got reordered in the source feed. Column A ("id"/row index) stays put; every
# other field (B:AD - match id, teams, scores, odds, ...) needs to move to the
# row where that fixture's data now belongs.
#
# This swaps/rotates the B:AD payload between affected rows while leaving
# column A (and the row's own formatting) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($sheet, $rowA, $rowB) {
    $rangeA = $sheet.Range("B$rowA`:AD$rowA")
    $rangeB = $sheet.Range("B$rowB`:AD$rowB")
    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2
    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

function Rotate-RowData($sheet, $rows) {
    # Shifts B:AD payloads through the given rows in a cycle:
    # rows[0] <- rows[1] <- rows[2] <- ... <- rows[0]
    $ranges = @()
    $originalValues = @()
    foreach ($r in $rows) {
        $rng = $sheet.Range("B$r`:AD$r")
        $ranges += $rng
        $originalValues += $rng.Value2
    }
    $count = $rows.Count
    for ($i = 0; $i -lt $count; $i++) {
        $srcIndex = ($i + 1) % $count
        $ranges[$i].Value2 = $originalValues[$srcIndex]
    }
}

# Two fixtures on 2024-03-?? (row 228/229) had their data swapped.
Swap-RowData $ws 228 229

# Two fixtures on the next match-day (row 252/253) swapped.
Swap-RowData $ws 252 253

# Two fixtures (row 305/306) swapped.
Swap-RowData $ws 305 306

# Two fixtures (row 312/313) swapped.
Swap-RowData $ws 312 313

# Three fixtures on the same match-day (row 338/339/340) rotated.
Rotate-RowData $ws @(338, 339, 340)
